$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 (columns P and Q), copying the header
# formatting (bold/centered/bordered style) from the existing O1 header cell.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14

$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1").Value = 15

$excel.CutCopyMode = 0

# For every data row (2-25), swap the I/K and M/O values and append the
# two new data columns P and Q (both value 2, default formatting).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2 (new column)
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2 (new column)
}
